$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cgn"
$ws.Range("C2").Value = "F11r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.0726045
$ws.Range("H2").Value = 0.145209
$ws.Range("I2").Value = 0.6461094194719614
$ws.Range("J2").Value = 0.5489714983497851
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 33.8141135
$ws.Range("N2").Value = 67.62822700000001
$ws.Range("O2").Value = 0.8395820055673167
$ws.Range("P2").Value = 0.7840733132753086
$ws.Range("Q2").Value = 2.45505680361075
$ws.Range("R2").Value = 9.820227214443001
$ws.Range("S2").Value = 0.5424618422162041
$ws.Range("T2").Value = 0.4304339016048266

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cgn"
$ws.Range("C3").Value = "F11r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.0726045
$ws.Range("H3").Value = 0.145209
$ws.Range("I3").Value = 0.6461094194719614
$ws.Range("J3").Value = 0.5489714983497851
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.08656933333333333
$ws.Range("N3").Value = 0.259708
$ws.Range("O3").Value = 0.002149459115662623
$ws.Range("P3").Value = 0.003011022483911989
$ws.Range("Q3").Value = 0.006285323162
$ws.Range("R3").Value = 0.037711938972
$ws.Range("S3").Value = 0.001388785781399493
$ws.Range("T3").Value = 0.001652965524558056

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cgn"
$ws.Range("C4").Value = "F11r"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.0726045
$ws.Range("H4").Value = 0.145209
$ws.Range("I4").Value = 0.6461094194719614
$ws.Range("J4").Value = 0.5489714983497851
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.533859666666667
$ws.Range("N4").Value = 7.601579
$ws.Range("O4").Value = 0.06291405453424448
$ws.Range("P4").Value = 0.088131768302221
$ws.Range("Q4").Value = 0.1839696141685
$ws.Range("R4").Value = 1.103817685011
$ws.Range("S4").Value = 0.04064936325174803
$ws.Range("T4").Value = 0.04838182889708636

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Cgn"
$ws.Range("C5").Value = "F11r"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.0726045
$ws.Range("H5").Value = 0.145209
$ws.Range("I5").Value = 0.6461094194719614
$ws.Range("J5").Value = 0.5489714983497851
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7582734999999999
$ws.Range("N5").Value = 1.516547
$ws.Range("O5").Value = 0.01882742795840407
$ws.Range("P5").Value = 0.0175826586586061
$ws.Range("Q5").Value = 0.05505406833074999
$ws.Range("R5").Value = 0.220216273323
$ws.Range("S5").Value = 0.01216457854835463
$ws.Range("T5").Value = 0.009652378468787816

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Cgn"
$ws.Range("C6").Value = "F11r"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.0726045
$ws.Range("H6").Value = 0.145209
$ws.Range("I6").Value = 0.6461094194719614
$ws.Range("J6").Value = 0.5489714983497851
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.9893793333333333
$ws.Range("N6").Value = 2.968138
$ws.Range("O6").Value = 0.02456563248203607
$ws.Range("P6").Value = 0.03441222547381506
$ws.Range("Q6").Value = 0.071833391807
$ws.Range("R6").Value = 0.4310003508419999
$ws.Range("S6").Value = 0.01587208654192989
$ws.Range("T6").Value = 0.0188913309799109

# Row 7
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Cgn"
$ws.Range("C7").Value = "F11r"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.0726045
$ws.Range("H7").Value = 0.145209
$ws.Range("I7").Value = 0.6461094194719614
$ws.Range("J7").Value = 0.5489714983497851
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.092743
$ws.Range("N7").Value = 6.278229
$ws.Range("O7").Value = 0.05196142034233613
$ws.Range("P7").Value = 0.0727890118061372
$ws.Range("Q7").Value = 0.1519425591435
$ws.Range("R7").Value = 0.911655354861
$ws.Range("S7").Value = 0.03357276313232536
$ws.Range("T7").Value = 0.03995909287461533

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Cgn"
$ws.Range("C8").Value = "F11r"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.03976733333333333
$ws.Range("H8").Value = 0.119302
$ws.Range("I8").Value = 0.3538905805280386
$ws.Range("J8").Value = 0.4510285016502149
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 33.8141135
$ws.Range("N8").Value = 67.62822700000001
$ws.Range("O8").Value = 0.8395820055673167
$ws.Range("P8").Value = 0.7840733132753086
$ws.Range("Q8").Value = 1.344697122925667
$ws.Range("R8").Value = 8.068182737554
$ws.Range("S8").Value = 0.2971201633511126
$ws.Range("T8").Value = 0.353639411670482

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Cgn"
$ws.Range("C9").Value = "F11r"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.03976733333333333
$ws.Range("H9").Value = 0.119302
$ws.Range("I9").Value = 0.3538905805280386
$ws.Range("J9").Value = 0.4510285016502149
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.08656933333333333
$ws.Range("N9").Value = 0.259708
$ws.Range("O9").Value = 0.002149459115662623
$ws.Range("P9").Value = 0.003011022483911989
$ws.Range("Q9").Value = 0.003442631535111111
$ws.Range("R9").Value = 0.030983683816
$ws.Range("S9").Value = 0.0007606733342631299
$ws.Range("T9").Value = 0.001358056959353933

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Cgn"
$ws.Range("C10").Value = "F11r"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.03976733333333333
$ws.Range("H10").Value = 0.119302
$ws.Range("I10").Value = 0.3538905805280386
$ws.Range("J10").Value = 0.4510285016502149
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.533859666666667
$ws.Range("N10").Value = 7.601579
$ws.Range("O10").Value = 0.06291405453424448
$ws.Range("P10").Value = 0.088131768302221
$ws.Range("Q10").Value = 0.1007648419842222
$ws.Range("R10").Value = 0.9068835778579999
$ws.Range("S10").Value = 0.02226469128249646
$ws.Range("T10").Value = 0.03974993940513464

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Cgn"
$ws.Range("C11").Value = "F11r"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.03976733333333333
$ws.Range("H11").Value = 0.119302
$ws.Range("I11").Value = 0.3538905805280386
$ws.Range("J11").Value = 0.4510285016502149
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.7582734999999999
$ws.Range("N11").Value = 1.516547
$ws.Range("O11").Value = 0.01882742795840407
$ws.Range("P11").Value = 0.0175826586586061
$ws.Range("Q11").Value = 0.03015451503233333
$ws.Range("R11").Value = 0.180927090194
$ws.Range("S11").Value = 0.006662849410049441
$ws.Range("T11").Value = 0.00793028018981829

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Cgn"
$ws.Range("C12").Value = "F11r"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.03976733333333333
$ws.Range("H12").Value = 0.119302
$ws.Range("I12").Value = 0.3538905805280386
$ws.Range("J12").Value = 0.4510285016502149
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.9893793333333333
$ws.Range("N12").Value = 2.968138
$ws.Range("O12").Value = 0.02456563248203607
$ws.Range("P12").Value = 0.03441222547381506
$ws.Range("Q12").Value = 0.03934497774177777
$ws.Range("R12").Value = 0.3541047996759999
$ws.Range("S12").Value = 0.008693545940106188
$ws.Range("T12").Value = 0.01552089449390416

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Cgn"
$ws.Range("C13").Value = "F11r"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.03976733333333333
$ws.Range("H13").Value = 0.119302
$ws.Range("I13").Value = 0.3538905805280386
$ws.Range("J13").Value = 0.4510285016502149
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.092743
$ws.Range("N13").Value = 6.278229
$ws.Range("O13").Value = 0.05196142034233613
$ws.Range("P13").Value = 0.0727890118061372
$ws.Range("Q13").Value = 0.08322280846199999
$ws.Range("R13").Value = 0.7490052761579999
$ws.Range("S13").Value = 0.01838865721001076
$ws.Range("T13").Value = 0.03282991893152187
